$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ B=0.5528496190351575; C=0.08710290462138914; D=0.113681439038146; E=0.1308403623878931; F=2.013411602932891; I=1.327276427429808; J=0.1826122813970965; K=0.33300711553062; L=0.2848554801963203; M=0.1863783687266398; O=5.381132015499006 }
    3 = @{ B=0.5198161574487301; C=0.08495753218657143; D=0.11230146106967; E=0.1312754649252188; F=2.023880063854193; I=1.337344355234844; J=0.1837058974389549; K=0.3011939872233; L=0.2829682393676904; M=0.1803741292745684; O=5.416023019337288 }
    4 = @{ B=0.4996749214235479; C=0.08362541347447916; D=0.1114903008318748; E=0.1315773510262837; F=2.031125080598912; I=1.344034547124426; J=0.1844162566822192; K=0.2816919988135567; L=0.2819061833775365; M=0.1767492342256602; O=5.439667239191238 }
    5 = @{ B=0.491503434854252; C=0.08307885008526483; D=0.1111688939072266; E=0.1317091284700993; F=2.034283263114041; I=1.34688882031854; J=0.1847155255152977; K=0.2737532130162634; L=0.2814978089927109; M=0.1752877284184038; O=5.449861444732889 }
    6 = @{ B=0.4901487735340027; C=0.08298786992236984; D=0.1111160786882053; E=0.1317315395535665; F=2.034820112912179; I=1.347370502356871; J=0.1847658107509602; K=0.2724355088290622; L=0.281431476760325; M=0.1750459971017975; O=5.451587961510526 }
    7 = @{ B=0.4995645703208993; C=0.08361805734297434; D=0.1114859291125043; E=0.1315790927379972; F=2.031166839424877; I=1.344072522572535; J=0.184420253058744; K=0.2815848987152805; L=0.2819005768784422; M=0.1767294602136147; O=5.439802457799004 }
    8 = @{ B=0.541430711753577; C=0.0863662741764486; D=0.1131981527465555; E=0.1309831909213131; F=2.016851639868598; I=1.33064238977677; J=0.1829813016901563; K=0.322031743314767; L=0.2841847456799158; M=0.184295379333328; O=5.392701765992342 }
    9 = @{ B=0.62462757568548; C=0.09163701885213982; D=0.1168404206387095; E=0.1300892273501848; F=1.995254335404766; I=1.308335210573766; J=0.1804672160897161; K=0.4015774117172839; L=0.289427388168022; M=0.1996164312810436; O=5.31793894769811 }
    10 = @{ B=0.686394457809314; C=0.09543664501090632; D=0.119687244047789; E=0.1295985569414153; F=1.983320862852835; I=1.294395729772067; J=0.1788066158602826; K=0.4601389985058972; L=0.2937398086165643; M=0.2111619338697182; O=5.273714109866233 }
    11 = @{ B=0.7146279456864875; C=0.09714929256796268; D=0.1210188971619885; E=0.1294111435265251; F=1.978743561948257; I=1.288584795770532; J=0.1780914321352078; K=0.4868021117776209; L=0.295800705422181; M=0.2164758458219751; O=5.255913687310169 }
    12 = @{ B=0.7253381207084431; C=0.09779553660770546; D=0.1215283726899372; E=0.1293452998323072; F=1.977132440150925; I=1.286460482489957; J=0.1778263792147285; K=0.4969016026969939; L=0.2965952782313792; M=0.2184968436715025; O=5.249505954808626 }
    13 = @{ B=0.7230306695227569; C=0.09765645898322362; D=0.1214184173472432; E=0.1293592528260596; F=1.977473992442278; I=1.28691460536507; J=0.1778832066246148; K=0.4947263836164382; L=0.2964235248020941; M=0.2180611995747981; O=5.250871174795321 }
    14 = @{ B=0.7155087061836696; C=0.09720250564932087; D=0.121060707990523; E=0.12940562392291; F=1.978608566081213; I=1.288408501318642; J=0.178069510477342; K=0.4876329512113955; L=0.2958657922861931; M=0.2166419403486373; O=5.255379849628667 }
    15 = @{ B=0.7109037093402719; C=0.09692414611262734; D=0.1208422771070303; E=0.1294346944364939; F=1.97931943328873; I=1.289333471760195; J=0.1781843782014425; K=0.4832883597727005; L=0.2955260054955744; M=0.2157737360817791; O=5.258184885502288 }
    16 = @{ B=0.6845519975013588; C=0.09532439910282164; D=0.1196009492593078; E=0.1296115232897534; F=1.983637112304045; I=1.294786150083812; J=0.1788541631946678; K=0.4583969175754987; L=0.2936071106552447; M=0.2108158878356363; O=5.274924014615749 }
    17 = @{ B=0.6684202568856392; C=0.09433893726043152; D=0.1188487722461957; E=0.1297291546489827; F=1.986503754494485; I=1.29826693183481; J=0.1792753491491661; K=0.4431323405392789; L=0.2924552538720278; M=0.2077901294716682; O=5.285786271864879 }
    18 = @{ B=0.6591544952548247; C=0.09377063787377438; D=0.1184195901467717; E=0.1298001835782578; F=1.988232711934394; I=1.300318895761535; J=0.1795213917706286; K=0.434354773542907; L=0.2918020753818098; M=0.2060556197747587; O=5.292252128308206 }
    19 = @{ B=0.6560194856737382; C=0.09357796664546925; D=0.1182748708218; E=0.129824812295487; F=1.988831877933521; I=1.301022230131078; J=0.1796053485158193; K=0.4313832381508291; L=0.2915825274861277; M=0.2054693499576032; O=5.294478839225633 }
    20 = @{ B=0.6701361897424363; C=0.09444399550720561; D=0.1189284861121536; E=0.1297162839143553; F=1.986190302820063; I=1.297891231735399; J=0.1792301212824192; K=0.4447570558964458; L=0.2925769051777465; M=0.2081116247030295; O=5.284607388649761 }
    21 = @{ B=0.7177175854833706; C=0.09733590531205039; D=0.12116563503713; E=0.1293918646694099; F=1.978271999448182; I=1.28796764151717; J=0.17801463198718; K=0.4897163942836187; L=0.2960292284743105; M=0.2170585753449572; O=5.25404651070528 }
    22 = @{ B=0.7489238328770398; C=0.09921252783986745; D=0.1226580668549673; E=0.1292097044650262; F=1.973809122808035; I=1.281925890437321; J=0.1772538747357553; K=0.5191156730865032; L=0.2983679853815886; M=0.2229567667280179; O=5.236013438279343 }
    23 = @{ B=0.7322587263010405; C=0.09821217356545731; D=0.1218587723580953; E=0.1293042011182219; F=1.976125948625906; I=1.285109897716282; J=0.1776568320400616; K=0.5034234749898019; L=0.297112235156419; M=0.2198041893446074; O=5.245460622178911 }
    24 = @{ B=0.6693603896854938; C=0.09439650410530476; D=0.1188924373325264; E=0.1297220921779445; F=1.986331762349835; I=1.298060927454866; J=0.1792505566603602; K=0.4440225275169212; L=0.2925218784588637; M=0.207966261001495; O=5.285139673319804 }
    25 = @{ B=0.6020059146137555; C=0.09022388613966115; D=0.1158248966158268; E=0.1303018029128573; F=2.000405168270852; I=1.313939231779145; J=0.1811145142223971; K=0.3800357570525534; L=0.2879278779314873; M=0.1954204662692725; O=5.336282699052845 }
}

foreach ($row in $data.Keys) {
    foreach ($col in $data[$row].Keys) {
        $ws.Range("$col$row").Value = $data[$row][$col]
    }
}
